# Update market-price / profit figures across multiple leve sheets
# (scheduled-runner refresh of cached Universalis price data).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 564.2963
$ws.Range("I107").Value = 574.35297
$ws.Range("J107").Value = 547.2
$ws.Range("K107").Value = 574.35297
$ws.Range("L107").Value = 547.2
$ws.Range("M107").Value = 1345.64703
$ws.Range("N107").Value = -4387.2

$ws.Range("H110").Value = 36208.332
$ws.Range("J110").Value = 36208.332
$ws.Range("L110").Value = 36208.332
$ws.Range("N110").Value = -44388.332

$ws.Range("H132").Value = 12289
$ws.Range("I132").Value = 3577.6667
$ws.Range("J132").Value = 51490
$ws.Range("K132").Value = 10733.0001
$ws.Range("L132").Value = 154470
$ws.Range("M132").Value = -8203.000100000001
$ws.Range("N132").Value = -159530

$ws.Range("H135").Value = 111112420
$ws.Range("I135").Value = 1157.6666
$ws.Range("J135").Value = 333334940
$ws.Range("K135").Value = 10418.9994
$ws.Range("L135").Value = 3000014460
$ws.Range("M135").Value = -7883.999400000001
$ws.Range("N135").Value = -3000019530

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1906.4529
$ws.Range("I61").Value = 1837.85
$ws.Range("J61").Value = 2117.5386
$ws.Range("K61").Value = 1837.85
$ws.Range("L61").Value = 2117.5386
$ws.Range("M61").Value = -1625.85
$ws.Range("N61").Value = -2541.5386

$ws.Range("H110").Value = 2769.45
$ws.Range("I110").Value = 2577.1667
$ws.Range("J110").Value = 4500
$ws.Range("K110").Value = 2577.1667
$ws.Range("L110").Value = 4500
$ws.Range("M110").Value = -532.1667000000002
$ws.Range("N110").Value = -8590

$ws.Range("H122").Value = 2112.5
$ws.Range("I122").Value = 1898.6
$ws.Range("J122").Value = 2469
$ws.Range("K122").Value = 5695.799999999999
$ws.Range("L122").Value = 7407
$ws.Range("M122").Value = -3245.799999999999
$ws.Range("N122").Value = -12307

$ws.Range("H136").Value = 1906.4529
$ws.Range("I136").Value = 1837.85
$ws.Range("J136").Value = 2117.5386
$ws.Range("K136").Value = 5513.549999999999
$ws.Range("L136").Value = 6352.6158
$ws.Range("M136").Value = -2963.549999999999
$ws.Range("N136").Value = -11452.6158

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 773.11536
$ws.Range("I80").Value = 1268.9166
$ws.Range("K80").Value = 1268.9166
$ws.Range("M80").Value = -270.9166

$ws.Range("H83").Value = 773.11536
$ws.Range("I83").Value = 1268.9166
$ws.Range("K83").Value = 6344.583000000001
$ws.Range("M83").Value = -1352.583000000001

$ws.Range("H99").Value = 2203.5715
$ws.Range("I99").Value = 1670
$ws.Range("K99").Value = 1670
$ws.Range("M99").Value = -172

$ws.Range("H107").Value = 1635.875
$ws.Range("I107").Value = 1578.2727
$ws.Range("J107").Value = 1762.6
$ws.Range("K107").Value = 1578.2727
$ws.Range("L107").Value = 1762.6
$ws.Range("M107").Value = 341.7273
$ws.Range("N107").Value = -5602.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 3074.9167
$ws.Range("I4").Value = 300
$ws.Range("J4").Value = 3327.182
$ws.Range("K4").Value = 300
$ws.Range("L4").Value = 3327.182
$ws.Range("M4").Value = -188
$ws.Range("N4").Value = -3551.182

$ws.Range("H100").Value = 30775.77
$ws.Range("J100").Value = 30775.77
$ws.Range("L100").Value = 30775.77
$ws.Range("N100").Value = -32939.77

$ws.Range("H119").Value = 30181.818
$ws.Range("J119").Value = 30181.818
$ws.Range("L119").Value = 30181.818
$ws.Range("N119").Value = -39857.818

$ws.Range("H140").Value = 38288.184
$ws.Range("J140").Value = 38288.184
$ws.Range("L140").Value = 38288.184
$ws.Range("N140").Value = -48648.184

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1034.5769
$ws.Range("I4").Value = 349.75
$ws.Range("J4").Value = 1159.091
$ws.Range("K4").Value = 1049.25
$ws.Range("L4").Value = 3477.273
$ws.Range("M4").Value = -937.25
$ws.Range("N4").Value = -3701.273

$ws.Range("H6").Value = 1217.3334
$ws.Range("J6").Value = 350
$ws.Range("L6").Value = 1050
$ws.Range("N6").Value = -1276

$ws.Range("H33").Value = 82.5
$ws.Range("J33").Value = 90
$ws.Range("L33").Value = 540
$ws.Range("N33").Value = -1106

$ws.Range("H46").Value = 884569.0600000001
$ws.Range("J46").Value = 1002418.25
$ws.Range("L46").Value = 3007254.75
$ws.Range("N46").Value = -3007436.75

$ws.Range("H61").Value = 892.5
$ws.Range("J61").Value = 890
$ws.Range("L61").Value = 2670
$ws.Range("N61").Value = -3100

$ws.Range("H113").Value = 1515668.6
$ws.Range("I113").Value = 2755307.8
$ws.Range("J113").Value = 554.1111
$ws.Range("K113").Value = 8265923.399999999
$ws.Range("L113").Value = 1662.3333
$ws.Range("M113").Value = -8263753.399999999
$ws.Range("N113").Value = -6002.3333

$ws.Range("H131").Value = 915.65
$ws.Range("I131").Value = 622.25
$ws.Range("J131").Value = 927.875
$ws.Range("K131").Value = 1866.75
$ws.Range("L131").Value = 2783.625
$ws.Range("M131").Value = 3173.25
$ws.Range("N131").Value = -12863.625

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 11600
$ws.Range("I5").Value = 1400
$ws.Range("J5").Value = 15000
$ws.Range("K5").Value = 1400
$ws.Range("L5").Value = 15000
$ws.Range("M5").Value = -1288
$ws.Range("N5").Value = -15224

$ws.Range("H102").Value = 41667268
$ws.Range("I102").Value = 83333336
$ws.Range("J102").Value = 1200
$ws.Range("K102").Value = 83333336
$ws.Range("L102").Value = 1200
$ws.Range("M102").Value = -83331714
$ws.Range("N102").Value = -4444

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H112").Value = 31385.6
$ws.Range("J112").Value = 31385.6
$ws.Range("L112").Value = 31385.6
$ws.Range("N112").Value = -34339.6

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H104").Value = 23388.572
$ws.Range("J104").Value = 23388.572
$ws.Range("L104").Value = 23388.572
$ws.Range("N104").Value = -30376.572

$ws.Range("H126").Value = 1428.3125
$ws.Range("I126").Value = 1350.5
$ws.Range("J126").Value = 1506.125
$ws.Range("K126").Value = 4051.5
$ws.Range("L126").Value = 4518.375
$ws.Range("M126").Value = -1581.5
$ws.Range("N126").Value = -9458.375

$ws.Range("H132").Value = 1437452.2
$ws.Range("I132").Value = 2030991.1
$ws.Range("J132").Value = 3066.4167
$ws.Range("K132").Value = 6092973.300000001
$ws.Range("L132").Value = 9199.250100000001
$ws.Range("M132").Value = -6090443.300000001
$ws.Range("N132").Value = -14259.2501

Write-Host "Updated 164 cells across 8 sheets"
